$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# -----------------------------------------------------------------
# 1) Update the "fecha de inicio" paragraph: change the date and add
#    the new justification text (why 31/12/2016 was chosen).
# -----------------------------------------------------------------
$idxFecha = Get-ParagraphIndexByText $d "Elegimos como fecha"
$pFecha = $d.Paragraphs.Item($idxFecha)
$pFecha.Range.Find.Execute(
    "01/01/2015 ya que todas las publicaciones pertenecen al año 2015",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "31/12/2016 ya que la última operación realizada de los datos migrados fue una compra el 30/12/2016",
    2
) | Out-Null

# -----------------------------------------------------------------
# 2) Insert a brand-new list paragraph right after it, explaining why
#    the migrated publications will all be finished by that date.
# -----------------------------------------------------------------
$pFecha.Range.InsertParagraphAfter() | Out-Null
$pNueva = $d.Paragraphs.Item($idxFecha + 1)
$pNueva.Range.InsertBefore(
    "Todas las publicaciones migradas estarán finalizadas ya que sus vencimientos son menores a esta fecha. Si bien hay compras/facturas con fechas sobre publicaciones que vencieron, nada nos asegura que estos datos sean erróneos ya que no conocemos el funcionamiento del sistema anterior."
) | Out-Null

# -----------------------------------------------------------------
# 3) Move the <w:lastRenderedPageBreak/> marker: it used to sit on the
#    "Las calificaciones migradas..." paragraph, and now belongs on the
#    "Cada vez que inicia el sistema..." paragraph instead.
# -----------------------------------------------------------------
$idxCadaVez = Get-ParagraphIndexByText $d "Cada vez que inicia el sistema"
$pCadaVez = $d.Paragraphs.Item($idxCadaVez)
$rCadaVez = $pCadaVez.Range.Duplicate
$rCadaVez.Collapse(1)
$pkgCadaVez = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00DF0290" w:rsidRDefault="00DF0290" w:rsidP="00DD5FFA"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Cada vez que inicia el sistema se hace un barrido sobre las publicaciones y se finaliza aquellas que hayan vencido. En caso de que sean subastas se generara la facturación correspondiente y se adjudicara la compra al ultimo usuario que oferto.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rCadaVez.InsertXML($pkgCadaVez) | Out-Null

$idxCalif = Get-ParagraphIndexByText $d "Las calificaciones migradas"
$pCalif = $d.Paragraphs.Item($idxCalif)
$rCalif = $pCalif.Range.Duplicate
$rCalif.Collapse(1)
$pkgCalif = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00766F24" w:rsidRDefault="00766F24" w:rsidP="00DD5FFA"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Las calificaciones migradas fueron migradas con la mitad de su valor ya que el rango de calificaciones era de 1 a 10 y ahora es de 1 a 5. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rCalif.InsertXML($pkgCalif) | Out-Null
